$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 38: add upstream_habitat_length_m=1000, species_codes changes to "RB, CT, CO, BT"
# (set first so new shared string "RB, CT, CO, BT" is registered before "RB, CO")
$ws.Range("K38").Value = 1000
$ws.Range("L38").Value = "RB, CT, CO, BT"

# Row 36: add hab_value=moderate, upstream_habitat_length_m=780, species_codes="RB, CO"
$ws.Range("H36").Value = "moderate"
$ws.Range("K36").Value = 780
$ws.Range("L36").Value = "RB, CO"

# Row 42: survey_length_m changes 0 -> 180, uav_survey=yes, hab_value=moderate, upstream_habitat_length_m=4500, species_codes="RB, CT"
$ws.Range("D42").Value = 180
$ws.Range("G42").Value = "yes"
$ws.Range("H42").Value = "moderate"
$ws.Range("K42").Value = 4500
$ws.Range("L42").Value = "RB, CT"

# Update selection to J36
$ws.Range("J36").Select()
